# Update the "2. Data reporter" section of the SDG indicator metadata sheet
# to reflect the new reporting organization / contact / phone / website.
#
# Cell layout (sheet "Лист 1"):
#   B6  = Organization
#   B7  = Contact person(s)/Focal point
#   B8  = Contact person's email   (unchanged)
#   B9  = Contact person's phone
#   B10 = Organization website

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit the website first, then organization/contact/phone - this mirrors the
# order the strings were (re)typed in and keeps the shared-string table in
# the same append order as the source edit.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B6").Value  = "The National Statistical Committee (Department of Digital Development and Sustainable Development Statistics)"
$ws.Range("B7").Value  = "Mambetaliev T.A."
$ws.Range("B9").Value  = "(0312) 62 56 07"

# Leave the cursor/selection on the contact-person cell, matching the saved
# workbook's last active selection.
$ws.Range("B7").Select()
